$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity for row 3 (40mm 12v fan) from 1 to 2
$ws.Range("C3").Value = 2

# Add quantity for row 4 (5x2mm magnets) which previously had no value
$ws.Range("C4").Value = 1

# Update the selected/active cell to C3
$ws.Range("C3").Select()
